$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.43"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.09%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "36.28"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.17%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.028"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.38%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07860"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.66%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.127"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.42%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-0.92%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9225"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.46%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09533"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-3.93%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1859"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.37%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08851"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.05%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03613"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.40%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09922"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.31%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001433"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-2.97%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005708"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.09%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.470"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.32%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.143"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.67%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "9.62%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.82%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1338"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.90%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.187"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.20%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.36%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04572"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.65%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.21%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004784"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-8.75%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-6.87%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004755"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "75.23%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01850"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.94%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04725"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007802"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.37%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1384"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.21%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007737"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.24%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002284"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.42%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01117"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "10.67%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006368"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.70%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.30%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.27%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "51.68"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "45.38%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.001903"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-29.16%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002104"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.30%"
